$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("K2").Value = "2025-11-06T07:03:01.991828+00:00"
$ws1.Range("K3").Value = "2025-11-06T07:03:01.991864+00:00"
$ws1.Range("K4").Value = "2025-11-06T07:03:01.991886+00:00"
$ws1.Range("K5").Value = "2025-11-06T07:03:04.726184+00:00"
$ws1.Range("K6").Value = "2025-11-06T07:03:04.726268+00:00"
$ws1.Range("K7").Value = "2025-11-06T07:03:04.726295+00:00"
$ws1.Range("K8").Value = "2025-11-06T07:03:07.107912+00:00"
$ws1.Range("K9").Value = "2025-11-06T07:03:10.023490+00:00"
$ws1.Range("K10").Value = "2025-11-06T07:03:10.023519+00:00"
$ws1.Range("K11").Value = "2025-11-06T07:03:10.023538+00:00"
$ws1.Range("K12").Value = "2025-11-06T07:03:12.892711+00:00"
$ws1.Range("K13").Value = "2025-11-06T07:03:12.892738+00:00"
$ws1.Range("K14").Value = "2025-11-06T07:03:12.892756+00:00"
$ws1.Range("K15").Value = "2025-11-06T07:03:12.892773+00:00"
$ws1.Range("K16").Value = "2025-11-06T07:03:18.088455+00:00"
$ws1.Range("K17").Value = "2025-11-06T07:03:20.861505+00:00"
$ws1.Range("K18").Value = "2025-11-06T07:03:23.117517+00:00"
$ws1.Range("K19").Value = "2025-11-06T07:03:23.117546+00:00"
$ws1.Range("K20").Value = "2025-11-06T07:03:25.451066+00:00"
$ws1.Range("K21").Value = "2025-11-06T07:03:28.248153+00:00"
$ws1.Range("K22").Value = "2025-11-06T07:03:28.248182+00:00"
$ws1.Range("K23").Value = "2025-11-06T07:03:28.248200+00:00"
$ws1.Range("K24").Value = "2025-11-06T07:03:31.003838+00:00"
$ws1.Range("K25").Value = "2025-11-06T07:03:31.003867+00:00"
$ws1.Range("K26").Value = "2025-11-06T07:03:31.003885+00:00"
$ws1.Range("K27").Value = "2025-11-06T07:03:33.445499+00:00"
$ws1.Range("K28").Value = "2025-11-06T07:03:33.445536+00:00"
$ws1.Range("K29").Value = "2025-11-06T07:03:33.445557+00:00"
$ws1.Range("K30").Value = "2025-11-06T07:03:33.445575+00:00"
$ws1.Range("K31").Value = "2025-11-06T07:03:35.739318+00:00"
$ws1.Range("K32").Value = "2025-11-06T07:03:38.110110+00:00"
$ws1.Range("K33").Value = "2025-11-06T07:03:38.110145+00:00"
$ws1.Range("K34").Value = "2025-11-06T07:03:38.110165+00:00"
$ws1.Range("K35").Value = "2025-11-06T07:03:40.932906+00:00"
$ws1.Range("K36").Value = "2025-11-06T07:03:40.932935+00:00"
$ws1.Range("K37").Value = "2025-11-06T07:03:40.932952+00:00"
$ws1.Range("K38").Value = "2025-11-06T07:03:43.403327+00:00"
$ws1.Range("K39").Value = "2025-11-06T07:03:43.403357+00:00"
$ws1.Range("K40").Value = "2025-11-06T07:03:43.403378+00:00"
$ws1.Range("K41").Value = "2025-11-06T07:03:43.403393+00:00"
$ws1.Range("K42").Value = "2025-11-06T07:03:43.403408+00:00"
$ws1.Range("K43").Value = "2025-11-06T07:03:43.403424+00:00"
$ws1.Range("K44").Value = "2025-11-06T07:03:43.403438+00:00"
$ws1.Range("K45").Value = "2025-11-06T07:03:43.403452+00:00"
$ws1.Range("K46").Value = "2025-11-06T07:03:46.260731+00:00"
$ws1.Range("K47").Value = "2025-11-06T07:03:46.260761+00:00"
$ws1.Range("K48").Value = "2025-11-06T07:03:51.480274+00:00"
$ws1.Range("K49").Value = "2025-11-06T07:03:51.480303+00:00"
$ws1.Range("K50").Value = "2025-11-06T07:03:54.403582+00:00"
$ws1.Range("K51").Value = "2025-11-06T07:03:54.403611+00:00"

$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows("2:6").Delete()
